$wb = $excel.ActiveWorkbook

# --- Locate the existing "Gumag Flame Ruins" sheet (last sheet in the workbook) ---
$gumag = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Create the new "Requi Water Ruins" sheet as a copy of Gumag's template,     ---
# --- placed immediately after it, BEFORE Gumag itself gets its new Area 7-9 cols ---
$gumag.Copy($null, $gumag)
$requi = $wb.Worksheets.Item($wb.Worksheets.Count)
$requi.Name = "レクイの水遺跡 (Requi Water Ruins)"

# ============================================================
# Fill in the new Requi Water Ruins sheet
# ============================================================

# Row 1: English area headers (columns B..K -> Area 1..Area 10)
$requi.Range("B1").Value = "Requi Water Ruins, Area 1"
$requi.Range("C1").Value = "Requi Water Ruins, Area 2"
$requi.Range("D1").Value = "Requi Water Ruins, Area 3"
$requi.Range("E1").Value = "Requi Water Ruins, Area 4"
$requi.Range("F1").Value = "Requi Water Ruins, Area 5"
$requi.Range("G1").Value = "Requi Water Ruins, Area 6"
$requi.Range("H1").Value = "Requi Water Ruins, Area 7"
$requi.Range("I1").Value = "Requi Water Ruins, Area 8"
$requi.Range("J1").Value = "Requi Water Ruins, Area 9"
$requi.Range("K1").Value = "Requi Water Ruins, Area 10"

# Row 2: Japanese area headers (columns B..K -> Area 1..Area 10)
$requi.Range("B2").Value = "レクイの水遺跡　エリア１"
$requi.Range("C2").Value = "レクイの水遺跡　エリア２"
$requi.Range("D2").Value = "レクイの水遺跡　エリア３"
$requi.Range("E2").Value = "レクイの水遺跡　エリア４"
$requi.Range("F2").Value = "レクイの水遺跡　エリア５"
$requi.Range("G2").Value = "レクイの水遺跡　エリア６"
$requi.Range("H2").Value = "レクイの水遺跡　エリア７"
$requi.Range("I2").Value = "レクイの水遺跡　エリア８"
$requi.Range("J2").Value = "レクイの水遺跡　エリア９"
$requi.Range("K2").Value = "レクイの水遺跡　エリア１０"

# The copied template (rows 3-16) carries over Gumag's old Day7-area filenames in
# B/D/F/G/I/J of row 16 and has no data yet in row 9/row 10 B..F - clear/replace
# so the new sheet matches its own (different) set of example filenames.
$requi.Range("B16").ClearContents()
$requi.Range("D16").ClearContents()
$requi.Range("F16").ClearContents()
$requi.Range("G16").ClearContents()
$requi.Range("I16").ClearContents()
$requi.Range("J16").ClearContents()

# Row 9 (Day 6): two example filenames
$requi.Range("B9").Value = "Day 06\028_25800156_189addc.xml"
$requi.Range("E9").Value = "Day 06\034_25805740_189c3ac.xml"

# Row 16 (Tag Only?): example filenames
$requi.Range("C16").Value = "Day 06\030_25803036_189b91c.xml "
$requi.Range("H16").Value = "Day 06\038_25810972_189d81c.xml"
$requi.Range("K16").Value = "Day 06\043_25815772_189eadc.xml"

# ============================================================
# Extend the existing Gumag Flame Ruins sheet with Areas 7, 8 and 9
# ============================================================

# Row 1: English headers for the new areas
$gumag.Range("H1").Value = "Gumag Flame Ruins, Area 7"
$gumag.Range("I1").Value = "Gumag Flame Ruins, Area 8"
$gumag.Range("J1").Value = "Gumag Flame Ruins, Area 9"

# Row 2: Japanese headers for the new areas
$gumag.Range("H2").Value = "グマグの炎遺跡　エリア７"
$gumag.Range("I2").Value = "グマグの炎遺跡　エリア８"
$gumag.Range("J2").Value = "グマグの炎遺跡　エリア９"

# Row 10 (Day 7): newly documented example filenames
$gumag.Range("B10").Value = "Day 07\010_25817148_189f03c.xml"
$gumag.Range("C10").Value = "Day 07\25913948_18b6a5c.xml"
$gumag.Range("D10").Value = "Day 07\012_25821500_18a013c.xml"
$gumag.Range("E10").Value = "Day 07\25917580_18b788c.xml"
$gumag.Range("F10").Value = "Day 07\014_25825308_18a101c.xml"

# Row 16 (Tag Only?): B/D/F change to new filenames, H/I/J added for areas 7-9
$gumag.Range("B16").Value = "Day 07\25912876_18b662c.xml"
$gumag.Range("D16").Value = "Day 07\25915708_18b713c.xml"
$gumag.Range("F16").Value = "Day 07\25919804_18b813c.xml"
$gumag.Range("H16").Value = "Day 07\25828156_18a1b3c.xml"
$gumag.Range("I16").Value = "Day 07\25829916_18a221c.xml"
$gumag.Range("J16").Value = "Day 07\25832044_18a2a6c.xml"

# ============================================================
# Selection / active-tab bookkeeping so the new sheet ends up active
# ============================================================
$gumag.Range("F17").Select()
$requi.Activate()
$requi.Range("K16").Select()
